$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2356265.5
$ws.Range("J17").Value = 2421542.5
$ws.Range("L17").Value = 7264627.5
$ws.Range("N17").Value = -7264963.5
$ws.Range("H19").Value = 1898.6666
$ws.Range("I19").Value = 2393.818
$ws.Range("J19").Value = 537
$ws.Range("K19").Value = 2393.818
$ws.Range("L19").Value = 537
$ws.Range("M19").Value = -2218.818
$ws.Range("N19").Value = -887
$ws.Range("H41").Value = 463.04
$ws.Range("J41").Value = 401
$ws.Range("L41").Value = 401
$ws.Range("N41").Value = -1281
$ws.Range("H62").Value = 5702.1113
$ws.Range("I62").Value = 7089.3076
$ws.Range("K62").Value = 7089.3076
$ws.Range("M62").Value = -6465.3076
$ws.Range("H65").Value = 5702.1113
$ws.Range("I65").Value = 7089.3076
$ws.Range("K65").Value = 35446.538
$ws.Range("M65").Value = -32326.538
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = 0
$ws.Range("H86").Value = 3752.2727
$ws.Range("I86").Value = 3166.3333
$ws.Range("J86").Value = 3972
$ws.Range("K86").Value = 3166.3333
$ws.Range("L86").Value = 3972
$ws.Range("M86").Value = -2043.3333
$ws.Range("N86").Value = -6218
$ws.Range("H89").Value = 3752.2727
$ws.Range("I89").Value = 3166.3333
$ws.Range("J89").Value = 3972
$ws.Range("K89").Value = 15831.6665
$ws.Range("L89").Value = 19860
$ws.Range("M89").Value = -10215.6665
$ws.Range("N89").Value = -31092
$ws.Range("H111").Value = 1883.6875
$ws.Range("J111").Value = 3161
$ws.Range("L111").Value = 9483
$ws.Range("N111").Value = -15617
$ws.Range("H116").Value = 8743.166999999999
$ws.Range("I116").Value = 5144.3335
$ws.Range("K116").Value = 5144.3335
$ws.Range("M116").Value = -1702.3335
$ws.Range("H132").Value = 2285.7346
$ws.Range("I132").Value = 2186.5557
$ws.Range("K132").Value = 6559.6671
$ws.Range("M132").Value = -4029.6671
$ws.Range("H137").Value = 15744.392
$ws.Range("I137").Value = 7354.7
$ws.Range("K137").Value = 22064.1
$ws.Range("M137").Value = -19514.1
$ws.Range("H141").Value = 4042.2856
$ws.Range("I141").Value = 4824.75
$ws.Range("K141").Value = 14474.25
$ws.Range("M141").Value = -9294.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3972.9683
$ws.Range("I32").Value = 3972.9683
$ws.Range("K32").Value = 3972.9683
$ws.Range("M32").Value = -3685.9683
$ws.Range("H45").Value = 7141.45
$ws.Range("I45").Value = 7905
$ws.Range("K45").Value = 7905
$ws.Range("M45").Value = -7528
$ws.Range("H132").Value = 4847.4194
$ws.Range("I132").Value = 3178.96
$ws.Range("J132").Value = 11799.333
$ws.Range("K132").Value = 9536.880000000001
$ws.Range("L132").Value = 35397.999
$ws.Range("M132").Value = -7006.880000000001
$ws.Range("N132").Value = -40457.999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 8231
$ws.Range("I24").Value = 8231
$ws.Range("K24").Value = 8231
$ws.Range("M24").Value = -7996
$ws.Range("H80").Value = 302.0625
$ws.Range("I80").Value = 244.71428
$ws.Range("J80").Value = 346.66666
$ws.Range("K80").Value = 244.71428
$ws.Range("L80").Value = 346.66666
$ws.Range("M80").Value = 753.28572
$ws.Range("N80").Value = -2342.66666
$ws.Range("H83").Value = 302.0625
$ws.Range("I83").Value = 244.71428
$ws.Range("J83").Value = 346.66666
$ws.Range("K83").Value = 1223.5714
$ws.Range("L83").Value = 1733.3333
$ws.Range("M83").Value = 3768.4286
$ws.Range("N83").Value = -11717.3333
$ws.Range("H107").Value = 4347.1763
$ws.Range("J107").Value = 7550
$ws.Range("L107").Value = 7550
$ws.Range("N107").Value = -11390
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 61979.25
$ws.Range("J80").Value = 61979.25
$ws.Range("L80").Value = 61979.25
$ws.Range("N80").Value = -64225.25
$ws.Range("H83").Value = 61979.25
$ws.Range("J83").Value = 61979.25
$ws.Range("L83").Value = 185937.75
$ws.Range("N83").Value = -197169.75
$ws.Range("H86").Value = 3347.45
$ws.Range("I86").Value = 2507.6365
$ws.Range("J86").Value = 4373.8887
$ws.Range("K86").Value = 2507.6365
$ws.Range("L86").Value = 4373.8887
$ws.Range("M86").Value = -1384.6365
$ws.Range("N86").Value = -6619.8887
$ws.Range("H87").Value = 47513
$ws.Range("J87").Value = 47513
$ws.Range("L87").Value = 47513
$ws.Range("N87").Value = -49885
$ws.Range("H89").Value = 3347.45
$ws.Range("I89").Value = 2507.6365
$ws.Range("J89").Value = 4373.8887
$ws.Range("K89").Value = 12538.1825
$ws.Range("L89").Value = 21869.4435
$ws.Range("M89").Value = -6922.182500000001
$ws.Range("N89").Value = -33101.4435
$ws.Range("H90").Value = 47513
$ws.Range("J90").Value = 47513
$ws.Range("L90").Value = 142539
$ws.Range("N90").Value = -154395
$ws.Range("H122").Value = 2768.8572
$ws.Range("I122").Value = 1503.625
$ws.Range("J122").Value = 4455.8335
$ws.Range("K122").Value = 4510.875
$ws.Range("L122").Value = 13367.5005
$ws.Range("M122").Value = -2060.875
$ws.Range("N122").Value = -18267.5005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 595.8570999999999
$ws.Range("I86").Value = 492.33334
$ws.Range("K86").Value = 1477.00002
$ws.Range("M86").Value = -291.0000199999999
$ws.Range("H89").Value = 595.8570999999999
$ws.Range("I89").Value = 492.33334
$ws.Range("K89").Value = 4431.00006
$ws.Range("M89").Value = 1496.99994
$ws.Range("H98").Value = 1370.2
$ws.Range("I98").Value = 1239.8
$ws.Range("J98").Value = 1631
$ws.Range("K98").Value = 3719.4
$ws.Range("L98").Value = 4893
$ws.Range("M98").Value = -2221.4
$ws.Range("N98").Value = -7889
$ws.Range("H121").Value = 2454631
$ws.Range("H129").Value = 2080.3076
$ws.Range("I129").Value = 566
$ws.Range("J129").Value = 3378.2856
$ws.Range("K129").Value = 1698
$ws.Range("L129").Value = 10134.8568
$ws.Range("M129").Value = 3302
$ws.Range("N129").Value = -20134.8568
$ws.Range("H131").Value = 3285.7036
$ws.Range("I131").Value = 1383.7778
$ws.Range("J131").Value = 4236.6665
$ws.Range("K131").Value = 4151.3334
$ws.Range("L131").Value = 12709.9995
$ws.Range("M131").Value = 888.6665999999996
$ws.Range("N131").Value = -22789.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 35495
$ws.Range("J15").Value = 35495
$ws.Range("L15").Value = 35495
$ws.Range("N15").Value = -36071
$ws.Range("H33").Value = 36149.668
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H81").Value = 35495
$ws.Range("J81").Value = 35495
$ws.Range("L81").Value = 35495
$ws.Range("N81").Value = -37491
$ws.Range("H84").Value = 35495
$ws.Range("J84").Value = 35495
$ws.Range("L84").Value = 106485
$ws.Range("N84").Value = -116469
$ws.Range("H92").Value = 27520.818
$ws.Range("J92").Value = 24760.111
$ws.Range("L92").Value = 24760.111
$ws.Range("N92").Value = -28504.111
$ws.Range("H102").Value = 2623.3333
$ws.Range("J102").Value = 4000
$ws.Range("L102").Value = 4000
$ws.Range("N102").Value = -7244
$ws.Range("H132").Value = 27036.092
$ws.Range("I132").Value = 37900.668
$ws.Range("J132").Value = 13998.6
$ws.Range("K132").Value = 113702.004
$ws.Range("L132").Value = 41995.8
$ws.Range("M132").Value = -111172.004
$ws.Range("N132").Value = -47055.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4544.409
$ws.Range("I132").Value = 4395.263
$ws.Range("K132").Value = 13185.789
$ws.Range("M132").Value = -10655.789
$ws.Range("H136").Value = 5546.2666
$ws.Range("I136").Value = 5047.654
$ws.Range("K136").Value = 15142.962
$ws.Range("M136").Value = -12592.962

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 170013.83
$ws.Range("I132").Value = 275598.66
$ws.Range("J132").Value = 26720.143
$ws.Range("K132").Value = 826795.98
$ws.Range("L132").Value = 80160.429
$ws.Range("M132").Value = -824265.98
